# All test cases are now run in Chrome, so mark rows 2-6's Runmode column
# as "Y" instead of "N" (matching row 7, which was already "Y").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C2:C6").Value = "Y"

# Move the active selection to C8, just below the table.
$ws.Range("C8").Select()
